{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The \"Semaine 7\" section lists five team members (Karim, Michelle, Miguel,\n// Rapha\u00ebl, Thibaud), each as a bulleted paragraph of the form\n// \"<Name> : \" (bold name run, then a plain run containing \" : \" with\n// nothing after it yet). This change appends the week's accomplishment\n// text for each person as a new run at the end of their paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map \"<Name> : \" (paragraph's full text, nothing typed after the colon\n// yet) -> the sentence to append for that person this week.\nconst additions = {\n  \"Karim : \": \"d\u00e9veloppement de la connexion.\",\n  \"Michelle : \": \"d\u00e9veloppement de l'inscription.\",\n  \"Miguel : \": \"d\u00e9veloppement de la liste des projets.\",\n  \"Rapha\u00ebl : \": \"ajout de triggers et autres contraintes dans la base de donn\u00e9es + Cr\u00e9ation des mod\u00e8les dans l'application.\",\n  \"Thibaud : \": \"avance dans le module Project.\",\n};\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (Object.prototype.hasOwnProperty.call(additions, text)) {\n    paragraph.insertText(additions[text], Word.InsertLocation.end);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# The \"Semaine 7\" section lists five team members (Karim, Michelle, Miguel,\n# Rapha\u00ebl, Thibaud), each as a bulleted paragraph of the form\n# \"<Name> : \" (bold name run, then a plain run containing \" : \" with\n# nothing after it yet). This change appends the week's accomplishment\n# text for each person at the end of their paragraph (i.e. right before\n# the paragraph mark), adding it as a new run.\n\n$d = $word.ActiveDocument\n\n# Map \"<Name> : \" (paragraph text, without the trailing paragraph mark,\n# nothing typed after the colon yet) -> the sentence to append for that\n# person this week.\n$additions = @{\n    \"Karim : \"    = \"d\u00e9veloppement de la connexion.\"\n    \"Michelle : \" = \"d\u00e9veloppement de l'inscription.\"\n    \"Miguel : \"   = \"d\u00e9veloppement de la liste des projets.\"\n    \"Rapha\u00ebl : \"  = \"ajout de triggers et autres contraintes dans la base de donn\u00e9es + Cr\u00e9ation des mod\u00e8les dans l'application.\"\n    \"Thibaud : \"  = \"avance dans le module Project.\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13)\n    if ($additions.ContainsKey($text)) {\n        $p.Range.InsertAfter($additions[$text])\n    }\n}\n"}
